$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $oldStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $oldStyle
}

Set-TextValue $ws "D2" '29.316.23'
Set-TextValue $ws "E2" '  +0.00%  '
Set-TextValue $ws "D3" '1.843.17'
Set-TextValue $ws "D4" '0.9962'
Set-TextValue $ws "E4" '  -0.13%  '
Set-TextValue $ws "D5" '239.74'
Set-TextValue $ws "E5" '  -0.08%  '
Set-TextValue $ws "D6" '0.6269'
Set-TextValue $ws "E6" '  +0.11%  '
Set-TextValue $ws "D7" '0.9980'
Set-TextValue $ws "E7" '  +0.05%  '
Set-TextValue $ws "D8" '0.07492'
Set-TextValue $ws "E8" '  -1.48%  '
Set-TextValue $ws "D9" '0.2897'
Set-TextValue $ws "E9" '  -0.02%  '
Set-TextValue $ws "D10" '24.40'
Set-TextValue $ws "E10" '  -1.29%  '
Set-TextValue $ws "D11" '0.07724'
Set-TextValue $ws "E11" '  -0.02%  '
Set-TextValue $ws "D12" '1.843.37'
Set-TextValue $ws "E12" '  -2.38%  '
Set-TextValue $ws "D13" '4.981'
Set-TextValue $ws "E13" '  -0.75%  '
Set-TextValue $ws "D14" '0.6791'
Set-TextValue $ws "D15" '0.00001054'
Set-TextValue $ws "E15" '  +0.55%  '
Set-TextValue $ws "D16" '81.94'
Set-TextValue $ws "E16" '  -1.18%  '
Set-TextValue $ws "D17" '6.174'
Set-TextValue $ws "E17" '  +0.68%  '
Set-TextValue $ws "D18" '29.368.37'
Set-TextValue $ws "E18" '  +0.06%  '
Set-TextValue $ws "D19" '228.77'
Set-TextValue $ws "E19" '  +0.49%  '
Set-TextValue $ws "E20" '  -0.05%  '
Set-TextValue $ws "D21" '0.9978'
Set-TextValue $ws "E21" '  +0.06%  '
Set-TextValue $ws "D22" '7.481'
Set-TextValue $ws "E22" '  +0.39%  '
Set-TextValue $ws "D23" '0.9976'
Set-TextValue $ws "E23" '  -0.09%  '
Set-TextValue $ws "D24" '158.31'
Set-TextValue $ws "E24" '  -0.04%  '
Set-TextValue $ws "D25" '0.1372'
Set-TextValue $ws "E25" '  -0.78%  '
Set-TextValue $ws "D26" '8.417'
Set-TextValue $ws "E26" '  +0.09%  '
Set-TextValue $ws "D27" '17.49'
Set-TextValue $ws "E27" '  -0.84%  '
Set-TextValue $ws "D28" '0.06455'
Set-TextValue $ws "E28" '  +15.36%  '
Set-TextValue $ws "D29" '1.421'
Set-TextValue $ws "E29" '  +1.52%  '
Set-TextValue $ws "E30" '  +1.60%  '
Set-TextValue $ws "D31" '4.087'
Set-TextValue $ws "E31" '  -0.38%  '
Set-TextValue $ws "D32" '4.090'
Set-TextValue $ws "E32" '  +0.88%  '
Set-TextValue $ws "D33" '1.831'
Set-TextValue $ws "E33" '  +0.19%  '
Set-TextValue $ws "D34" '1.139'
Set-TextValue $ws "E34" '  -1.81%  '
Set-TextValue $ws "D35" '0.6968'
Set-TextValue $ws "E35" '  +0.08%  '
Set-TextValue $ws "D36" '2.574'
Set-TextValue $ws "E36" '  -0.27%  '
Set-TextValue $ws "D37" '1.268.85'
Set-TextValue $ws "E37" '  +3.70%  '
Set-TextValue $ws "E38" '  +4.33%  '
Set-TextValue $ws "D39" '0.01834'
Set-TextValue $ws "E39" '  +1.84%  '
Set-TextValue $ws "D40" '6.677'
Set-TextValue $ws "E40" '  +5.26%  '
Set-TextValue $ws "D41" '0.9140'
Set-TextValue $ws "E41" '  +1.47%  '
Set-TextValue $ws "D42" '0.9973'
Set-TextValue $ws "E42" '  +0.02%  '
Set-TextValue $ws "D43" '2.007.58'
Set-TextValue $ws "E43" '  -18.42%  '
Set-TextValue $ws "D44" '101.20'
Set-TextValue $ws "E44" '  +0.02%  '
Set-TextValue $ws "D45" '66.17'
Set-TextValue $ws "E45" '  +1.11%  '
Set-TextValue $ws "B46" 'Aptos'
Set-TextValue $ws "C46" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws "D46" '7.083'
Set-TextValue $ws "E46" '  -1.38%  '
Set-TextValue $ws "B47" 'RenderToken'
Set-TextValue $ws "C47" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws "D47" '1.721'
Set-TextValue $ws "E47" '  +2.65%  '
Set-TextValue $ws "D48" '0.1163'
Set-TextValue $ws "E48" '  +2.25%  '
Set-TextValue $ws "D49" '9.037'
Set-TextValue $ws "E49" '  +0.70%  '
Set-TextValue $ws "B50" 'BabyDogeCoin'
Set-TextValue $ws "C50" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws "D50" '0.00000000116'
Set-TextValue $ws "E50" '  -1.88%  '
Set-TextValue $ws "B51" 'TheSandbox'
Set-TextValue $ws "C51" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws "D51" '0.3961'
Set-TextValue $ws "E51" '  -0.65%  '
